# Generate Report for handback
# - Marks the in-flight handoff as handed back ("Handed back: in sync with en-US")
# - Records the "Latest Target File" / "Latest Handback File" (columns E/F) for
#   the two in-progress rows on each locale sheet, with hyperlinks matching the
#   style already used elsewhere in the sheet
# - Stamps the "Latest Handback DateTime" (column G) with the handback time

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$mdFile  = "a523183c-ca2a-404f-8c03-59116d38d675.md"
$mdUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/b48c5729623d1dfcbddce59abf29129c0422064b/e2e/a523183c-ca2a-404f-8c03-59116d38d675.md"

$zhXlf   = "a523183c-ca2a-404f-8c03-59116d38d675.ee9cdfa45bec56724d36bd8cdb71fb75082306a0.zh-cn.xlf"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f41eba0af36976709ce8d57a92905e7c07722739/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/a523183c-ca2a-404f-8c03-59116d38d675.ee9cdfa45bec56724d36bd8cdb71fb75082306a0.zh-cn.xlf"

$deXlf   = "a523183c-ca2a-404f-8c03-59116d38d675.ee9cdfa45bec56724d36bd8cdb71fb75082306a0.de-de.xlf"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3617d0ef9e698b7adab3f1736050adb0e2bbcfd5/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/a523183c-ca2a-404f-8c03-59116d38d675.ee9cdfa45bec56724d36bd8cdb71fb75082306a0.de-de.xlf"

$zhHandbackTime = "2016-01-26 07:40:36"
$deHandbackTime = "2016-01-26 07:40:55"

# ---- Overview sheet: status text is shared with the locale sheets, so it
# needs to be refreshed here too (rows for both tracked files) ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $statusText
$overview.Range("C2").Value = $statusText
$overview.Range("B3").Value = $statusText
$overview.Range("C3").Value = $statusText

# ---- Locale sheets: zh-cn and de-de share the same row layout ----
$locales = @(
  @{ Name = "zh-cn"; Xlf = $zhXlf; XlfUrl = $zhXlfUrl; HandbackTime = $zhHandbackTime },
  @{ Name = "de-de"; Xlf = $deXlf; XlfUrl = $deXlfUrl; HandbackTime = $deHandbackTime }
)

foreach ($locale in $locales) {
  $ws = $wb.Worksheets.Item($locale.Name)

  # Status column (B) for the two active rows
  $ws.Range("B2").Value = $statusText
  $ws.Range("B3").Value = $statusText

  # Latest Target File (E) / Latest Handback File (F) for rows 2 and 3
  $ws.Range("E2").Value = $mdFile
  $ws.Range("E2").Style = "Hyperlink"
  $ws.Hyperlinks.Add($ws.Range("E2"), $mdUrl, $null, $null, $mdFile)

  $ws.Range("F2").Value = $locale.Xlf
  $ws.Range("F2").Style = "Hyperlink"
  $ws.Hyperlinks.Add($ws.Range("F2"), $locale.XlfUrl, $null, $null, $locale.Xlf)

  $ws.Range("E3").Value = $mdFile
  $ws.Range("E3").Style = "Hyperlink"
  $ws.Hyperlinks.Add($ws.Range("E3"), $mdUrl, $null, $null, $mdFile)

  $ws.Range("F3").Value = $locale.Xlf
  $ws.Range("F3").Style = "Hyperlink"
  $ws.Hyperlinks.Add($ws.Range("F3"), $locale.XlfUrl, $null, $null, $locale.Xlf)

  # Latest Handback DateTime (G) for rows 2 and 3
  $ws.Range("G2").Value = $locale.HandbackTime
  $ws.Range("G3").Value = $locale.HandbackTime
}
